$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 holds an ID that looks numeric ("137463") but the source workbook stores
# it as text. A bare assignment gets auto-coerced to a Number by Excel's usual
# "looks like a number" input parsing, so force text with a leading
# apostrophe (the standard Excel quote-prefix trick), then strip the
# quote-prefix formatting flag back off so only the bold header style
# remains on the cell.
$ws.Range("A2").Value = "'137463"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = "COMPREHENDING THE TRANSPORT PROPERTIES OF PROTIC IONIC LIQUIDS USING NMR"
$ws.Range("C2").Value = "Poster Presentation"
$ws.Range("D2").Value = "Andrea Mele;"
$ws.Range("E2").Value = "Department of Chemistry, Materials and Chemical Engineering “Giulio Natta” / Politecnico di Milano"

$ws.Range("A2:E2").Font.Bold = $true
